$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the formula in D3 (row 3 previously had no cells at all)
$ws.Range("D3").Formula = "=COUNTA(D5:D20)"

# Move the student name from D11 to D18
$ws.Range("D18").Value = $ws.Range("D11").Value()
$ws.Range("D11").Value = ""

# Update row 18 height
$ws.Rows.Item(18).RowHeight = 29

# Update the active cell selection
$ws.Range("D3").Select()
